# More work on presentation.
#
# Merge the two runs of the "Mission: Help companies manage their
# templates" bullet on the "Templafy" slide into a single run, keeping
# the (dirty="0") run properties of the second run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$prefix = "Mission: Help "
$suffix = "companies manage their templates"

# Remove the leading run's text ("Mission: Help "). What remains is a
# single run - the one that already carries dirty="0" - that starts
# with "companies manage their templates".
$full = $tr.Text
$prefixStart = $full.IndexOf($prefix)
if ($prefixStart -ge 0) {
    $tr.Characters($prefixStart + 1, $prefix.Length).Text = ""
}

# Re-insert the removed text at the front of that same run so the
# whole sentence now lives in one run with one set of run properties.
$full2 = $tr.Text
$suffixStart = $full2.IndexOf($suffix)
if ($suffixStart -ge 0) {
    $tr.Characters($suffixStart + 1, $suffix.Length).Text = $prefix + $suffix
}
